$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the front (column A), shifting the existing
# "Headcount ... 0.5 or 1" table one column to the right (B:G).
$ws.Columns("A:A").Insert()

# Populate the new "Grad Level" column header and its value for row 2
$ws.Range("A1").Value = "Grad Level"
$ws.Range("A2").Value = "Ph.D. Student or M.S. Student"

# Size the new column A (matches the other bestFit-style header columns)
$ws.Columns("A:A").ColumnWidth = 27.166666666666668

# View changes: zoom to 85% and move the selection to B10
$excel.ActiveWindow.Zoom = 85
$ws.Range("B10").Select()
